$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 12120480
$ws.Range("I19").Value = 10435478
$ws.Range("K19").Value = 10435478
$ws.Range("M19").Value = -10435303

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3092.3076
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3092.3076
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3092.3076
$ws.Range("N74").Value = -4964.3076
$ws.Range("M74").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3092.3076
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3092.3076
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 15461.538
$ws.Range("N77").Value = -24821.538
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1514.6666
$ws.Range("I100").Value = 1361
$ws.Range("J100").Value = 1591.5
$ws.Range("K100").Value = 1361
$ws.Range("L100").Value = 1591.5
$ws.Range("M100").Value = -820
$ws.Range("N100").Value = -2673.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2352.32
$ws.Range("I132").Value = 2354.9092
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 7064.7276
$ws.Range("L132").Value = 6999.999899999999
$ws.Range("M132").Value = -4534.7276
$ws.Range("N132").Value = -12059.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3330.0256
$ws.Range("I138").Value = 722.2941
$ws.Range("J138").Value = 5345.091
$ws.Range("K138").Value = 2166.8823
$ws.Range("L138").Value = 16035.273
$ws.Range("M138").Value = 2973.1177
$ws.Range("N138").Value = -26315.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2100.027
$ws.Range("I61").Value = 1584.4062
$ws.Range("K61").Value = 1584.4062
$ws.Range("M61").Value = -1372.4062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 557.3158
$ws.Range("I97").Value = 522.4375
$ws.Range("J97").Value = 743.3333
$ws.Range("K97").Value = 522.4375
$ws.Range("L97").Value = 743.3333
$ws.Range("M97").Value = -26.4375
$ws.Range("N97").Value = -1735.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2100.027
$ws.Range("I136").Value = 1584.4062
$ws.Range("K136").Value = 4753.2186
$ws.Range("M136").Value = -2203.2186

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3260.1177
$ws.Range("I86").Value = 3109.3845
$ws.Range("J86").Value = 3750
$ws.Range("K86").Value = 3109.3845
$ws.Range("L86").Value = 3750
$ws.Range("M86").Value = -1986.3845
$ws.Range("N86").Value = -5996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3260.1177
$ws.Range("I89").Value = 3109.3845
$ws.Range("J89").Value = 3750
$ws.Range("K89").Value = 15546.9225
$ws.Range("L89").Value = 18750
$ws.Range("M89").Value = -9930.9225
$ws.Range("N89").Value = -29982

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 4400
$ws.Range("I97").Value = 4400
$ws.Range("K97").Value = 4400
$ws.Range("M97").Value = -3409

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 20000
$ws.Range("N88").Value = -20812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("N91").Value = -22808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3934.739
$ws.Range("I134").Value = 1153
$ws.Range("K134").Value = 3459
$ws.Range("M134").Value = -924

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 24285.715
$ws.Range("J88").Value = 24285.715
$ws.Range("L88").Value = 72857.145
$ws.Range("N88").Value = -73713.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 24285.715
$ws.Range("J91").Value = 24285.715
$ws.Range("L91").Value = 72857.145
$ws.Range("N91").Value = -75821.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 4065
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4065
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 12195
$ws.Range("N94").Value = -13547
$ws.Range("M94").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I129").Value = 1933.8
$ws.Range("J129").Value = 17577.334
$ws.Range("K129").Value = 5801.4
$ws.Range("L129").Value = 52732.00199999999
$ws.Range("M129").Value = -801.3999999999996
$ws.Range("N129").Value = -62732.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 858.0909
$ws.Range("J131").Value = 969.11536
$ws.Range("L131").Value = 2907.34608
$ws.Range("N131").Value = -12987.34608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3971.76
$ws.Range("I80").Value = 2771.0715
$ws.Range("J80").Value = 5499.909
$ws.Range("K80").Value = 2771.0715
$ws.Range("L80").Value = 5499.909
$ws.Range("M80").Value = -1773.0715
$ws.Range("N80").Value = -7495.909

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3971.76
$ws.Range("I83").Value = 2771.0715
$ws.Range("J83").Value = 5499.909
$ws.Range("K83").Value = 13855.3575
$ws.Range("L83").Value = 27499.545
$ws.Range("M83").Value = -8863.3575
$ws.Range("N83").Value = -37483.545

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9092944
$ws.Range("I122").Value = 12501760
$ws.Range("K122").Value = 37505280
$ws.Range("M122").Value = -37502830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3765.647
$ws.Range("I132").Value = 3430.2856
$ws.Range("K132").Value = 10290.8568
$ws.Range("M132").Value = -7760.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2401
$ws.Range("I82").Value = 1802
$ws.Range("K82").Value = 1802
$ws.Range("M82").Value = -1441

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2401
$ws.Range("I85").Value = 1802
$ws.Range("K85").Value = 1802
$ws.Range("M85").Value = -554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 44476.668
$ws.Range("J134").Value = 44476.668
$ws.Range("L134").Value = 44476.668
$ws.Range("N134").Value = -54616.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 34317
$ws.Range("J140").Value = 34317
$ws.Range("L140").Value = 34317
$ws.Range("N140").Value = -44677
